# Update Unit Test/XML Templates for onlineResources
# Adds new lookup rows (online resource link/protocol/description fields and
# contact detail fields) to the "Feuil1" sheet, rows 69-80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gmdContact = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

# Fill rows in the same order the original authoring session used so that
# new shared-string entries are allocated in the expected sequence:
# 69, 70, 74-80, then 71-73.

# Row 69
$ws.Range("A69").Value = "online_resource_link"
$ws.Range("B69").Value = "online_resource_link"
$ws.Range("C69").Value = $gmdContact

# Row 70
$ws.Range("A70").Value = "online_resource_protocol"
$ws.Range("B70").Value = "online_resource_protocol"
$ws.Range("C70").Value = $gmdContact

# Row 74
$ws.Range("A74").Value = "contact_phone"
$ws.Range("B74").Value = "contacts{}.phone"
$ws.Range("C74").Value = $gmdContact

# Row 75
$ws.Range("A75").Value = "contact_facsimile"
$ws.Range("B75").Value = "contacts{}.facsimile"
$ws.Range("C75").Value = $gmdContact

# Row 76
$ws.Range("A76").Value = "contact_address"
$ws.Range("B76").Value = "contacts{}.address"
$ws.Range("C76").Value = $gmdContact

# Row 77
$ws.Range("A77").Value = "contact_city"
$ws.Range("B77").Value = "contacts{}.city"
$ws.Range("C77").Value = $gmdContact

# Row 78
$ws.Range("A78").Value = "contact_administrative_area"
$ws.Range("B78").Value = "contacts{}.administrative_area"
$ws.Range("C78").Value = $gmdContact

# Row 79
$ws.Range("A79").Value = "contact_postalcode"
$ws.Range("B79").Value = "contacts{}.postalcode"
$ws.Range("C79").Value = $gmdContact

# Row 80
$ws.Range("A80").Value = "contact_country"
$ws.Range("B80").Value = "contacts{}.country"
$ws.Range("C80").Value = $gmdContact

# Row 71
$ws.Range("A71").Value = "online_resource_description"
$ws.Range("B71").Value = "online_resource_description_en"
$ws.Range("C71").Value = $gmdContact

# Row 72
$ws.Range("A72").Value = "online_resource_description_other_lang_locale"
$ws.Range("B72").Value = "online_resource_description_locale"
$ws.Range("C72").Value = $gmdContact

# Row 73
$ws.Range("A73").Value = "online_resource_description_other_lang"
$ws.Range("B73").Value = "online_resource_description_fr"
$ws.Range("C73").Value = $gmdContact

# Move the viewport/selection to mirror where the author ended up: scrolled
# down with the newly added block (rows 69-80) selected.
$excel.ActiveWindow.ScrollRow = 64
[void]$ws.Range("A69:XFD80").Select()
